$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header info updates
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number that must stay TEXT (not be coerced to a
# number). Prefixing with an apostrophe forces Excel's literal-entry parser
# to store it as text while keeping the existing cell style intact.
$ws.Range("B3").Value = "'2570314725427075"

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 11.02.2025"

# Row 6
$ws.Range("B6").Value = "14.02."
$ws.Range("C6").Value = "15.02."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 94236955"
$ws.Range("E6").Value = "86,77-"

# Row 7
$ws.Range("B7").Value = "16.02."
$ws.Range("C7").Value = "17.02."
$ws.Range("D7").Value = "RECHNUNG VODAFONE GMBH 81423328"
$ws.Range("E7").Value = "39,65-"

# Row 8
$ws.Range("B8").Value = "20.02."
$ws.Range("C8").Value = "21.02."
$ws.Range("D8").Value = "BEITRAG Allianz SE K-91393869"
$ws.Range("E8").Value = "54,49-"

# Row 9 - previously blank, now filled with a new transaction
$ws.Range("B9").Value = "21.02."
$ws.Range("C9").Value = "22.02."
$ws.Range("D9").Value = "KARTENZ./21.02 REWE RO"
$ws.Range("E9").Value = "111,51-"
# Match the amount-column formatting used by the other transaction rows
# (right-aligned, no wrap, default vertical alignment).
$ws.Range("E9").HorizontalAlignment = -4152
$ws.Range("E9").VerticalAlignment = -4107
$ws.Range("E9").WrapText = $false

# Row 10 - previously blank, now filled with a new transaction
$ws.Range("B10").Value = "24.02."
$ws.Range("C10").Value = "25.02."
$ws.Range("D10").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E10").Value = "25,35-"
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4107
$ws.Range("E10").WrapText = $false

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 28.02.2025"
$ws.Range("E12").Value = "317,77-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 07.03.2025"
